# Leave Card update — 12/22/2023 10:59 AM
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$conv = $wb.Worksheets.Item("CONVERTION")

# --- Sheet1 data edits ---

# A10 "2021" header: match the bold styling already used by the 2022/2023
# year headers (A20 / A33).
$ws1.Range("A10").Font.Bold = $true

# Remove the FL(5-0-0) / 5-day entry that was recorded on 12/11/2023 (row 32)
$ws1.Range("B32").ClearContents()
$ws1.Range("D32").ClearContents()

# Row 42 (4/17/2023): record the 1.25 VL/SL earned amount
$ws1.Range("C42").Value = 1.25

# Rows 44 and 45 (5/18/2023, 6/17/2023): record the 1.25 VL/SL earned amount
$ws1.Range("C44").Value = 1.25
$ws1.Range("C45").Value = 1.25

# Clear out the pre-filled future period dates in rows 46-83
for ($r = 46; $r -le 83; $r++) {
    $ws1.Cells.Item($r, 1).ClearContents()
}

# Scroll Sheet1's frozen/split pane up so row 23 is the first visible row
# below the split (was row 30).
$aw = $excel.ActiveWindow
$aw.SplitRow = 22

# --- Active sheet / selection ---
# CONVERTION becomes the active tab, with B12 selected.
$conv.Activate()
$conv.Range("B12").Select()
